$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings that must stay literal text
# (matches the source inline-string cells with no style). Force a Text
# number format while assigning, then clear the format so no style index
# lingers on the cell - this avoids Excel's automatic text->number coercion.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '37.655.45'
$ws.Range("E2").Value = '  +0.24%  '
Set-TextValue $ws.Range("D3") '2.105.88'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  +0.11%  '
Set-TextValue $ws.Range("D5") '235.63'
$ws.Range("E5").Value = '  +0.16%  '
Set-TextValue $ws.Range("D6") '0.626'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("E8").Value = '  +0.03%  '
Set-TextValue $ws.Range("D9") '0.391'
$ws.Range("E9").Value = '  +1.03%  '
Set-TextValue $ws.Range("D10") '0.0782'
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("E11").Value = '  +1.36%  '
Set-TextValue $ws.Range("D12") '2.418.40'
$ws.Range("E12").Value = '  +1.39%  '
Set-TextValue $ws.Range("D13") '14.60'
$ws.Range("E13").Value = '  +0.22%  '
Set-TextValue $ws.Range("D14") '21.34'
$ws.Range("E14").Value = '  +0.67%  '
Set-TextValue $ws.Range("D15") '0.790'
$ws.Range("E15").Value = '  +1.13%  '
Set-TextValue $ws.Range("D16") '5.23'
$ws.Range("E16").Value = '  +0.21%  '
Set-TextValue $ws.Range("D17") '2.111.77'
$ws.Range("E17").Value = '  +2.30%  '
Set-TextValue $ws.Range("D18") '37.677.15'
$ws.Range("E18").Value = '  -0.13%  '
Set-TextValue $ws.Range("D19") '6.18'
$ws.Range("E19").Value = '  -0.65%  '
Set-TextValue $ws.Range("D20") '69.86'
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  +0.95%  '
Set-TextValue $ws.Range("D22") '227.05'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +1.44%  '
Set-TextValue $ws.Range("D25") '2.42'
$ws.Range("E25").Value = '  -3.10%  '
Set-TextValue $ws.Range("D26") '169.34'
$ws.Range("E26").Value = '  +1.09%  '
Set-TextValue $ws.Range("D27") '8.96'
$ws.Range("E27").Value = '  -1.13%  '
$ws.Range("E28").Value = '  +3.00%  '
Set-TextValue $ws.Range("D29") '1.43'
$ws.Range("E29").Value = '  -4.86%  '
Set-TextValue $ws.Range("D30") '19.41'
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("E31").Value = '  -0.31%  '
Set-TextValue $ws.Range("D32") '4.65'
$ws.Range("E32").Value = '  +2.52%  '
Set-TextValue $ws.Range("D33") '0.0623'
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D34") '2.56'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D35") '4.60'
$ws.Range("E35").Value = '  -0.44%  '
Set-TextValue $ws.Range("D36") '3.55'
$ws.Range("E36").Value = '  +6.01%  '
Set-TextValue $ws.Range("D37") '1.78'
$ws.Range("E37").Value = '  +1.16%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -5.32%  '
$ws.Range("E40").Value = '  -0.34%  '
Set-TextValue $ws.Range("D41") '0.0959'
$ws.Range("E41").Value = '  -0.10%  '
Set-TextValue $ws.Range("D42") '97.78'
$ws.Range("E42").Value = '  +1.91%  '
Set-TextValue $ws.Range("D43") '1.483.07'
$ws.Range("E43").Value = '  +0.61%  '
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("E45").Value = '  -1.48%  '
Set-TextValue $ws.Range("D46") '4.16'
$ws.Range("E46").Value = '  -9.77%  '
$ws.Range("E47").Value = '  +1.11%  '
Set-TextValue $ws.Range("D48") '15.58'
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("E50").Value = '  +2.85%  '
Set-TextValue $ws.Range("D51") '2.303.98'
$ws.Range("E51").Value = '  +1.40%  '
